# Auto-generated edit script: update profit-calc columns (H..N) on 43 rows
# across all 8 sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 931.9474
$rowvals[0,1] = 923.94116
$rowvals[0,2] = 1000
$rowvals[0,3] = 923.94116
$rowvals[0,4] = 1000
$rowvals[0,5] = 574.05884
$rowvals[0,6] = -3996
$ws.Range("H98:N98").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 3508
$rowvals[0,1] = 3000
$rowvals[0,2] = 3554.182
$rowvals[0,3] = 3000
$rowvals[0,4] = 3554.182
$rowvals[0,5] = 254
$rowvals[0,6] = -10062.182
$ws.Range("H113:N113").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 931.9474
$rowvals[0,1] = 923.94116
$rowvals[0,2] = 1000
$rowvals[0,3] = 2771.82348
$rowvals[0,4] = 3000
$rowvals[0,5] = -321.82348
$rowvals[0,6] = -7900
$ws.Range("H122:N122").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 88237120
$rowvals[0,1] = 55558440
$rowvals[0,2] = 125000630
$rowvals[0,3] = 500025960
$rowvals[0,4] = 1125005670
$rowvals[0,5] = -500023425
$rowvals[0,6] = -1125010740
$ws.Range("H135:N135").Value = $rowvals

$ws = $wb.Worksheets.Item("ARM")
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 993.36
$rowvals[0,1] = 920.7619
$rowvals[0,2] = 1374.5
$rowvals[0,3] = 920.7619
$rowvals[0,4] = 1374.5
$rowvals[0,5] = -807.7619
$rowvals[0,6] = -1600.5
$ws.Range("H2:N2").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4960.2
$rowvals[0,1] = 4960.2
$rowvals[0,2] = 0
$rowvals[0,3] = 4960.2
$rowvals[0,4] = 0
$rowvals[0,5] = -4768.2
$rowvals[0,6] = $null
$ws.Range("H28:N28").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 1736.0938
$rowvals[0,1] = 1675.2084
$rowvals[0,2] = 1918.75
$rowvals[0,3] = 1675.2084
$rowvals[0,4] = 1918.75
$rowvals[0,5] = -1298.2084
$rowvals[0,6] = -2672.75
$ws.Range("H45:N45").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4960.2
$rowvals[0,1] = 4960.2
$rowvals[0,2] = 0
$rowvals[0,3] = 4960.2
$rowvals[0,4] = 0
$rowvals[0,5] = -1965.2
$rowvals[0,6] = $null
$ws.Range("H99:N99").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 1795.8125
$rowvals[0,1] = 1802.5385
$rowvals[0,2] = 1766.6666
$rowvals[0,3] = 1802.5385
$rowvals[0,4] = 1766.6666
$rowvals[0,5] = 242.4614999999999
$rowvals[0,6] = -5856.6666
$ws.Range("H110:N110").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 993.36
$rowvals[0,1] = 920.7619
$rowvals[0,2] = 1374.5
$rowvals[0,3] = 920.7619
$rowvals[0,4] = 1374.5
$rowvals[0,5] = 1373.2381
$rowvals[0,6] = -5962.5
$ws.Range("H116:N116").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 1717.7142
$rowvals[0,1] = 1490.375
$rowvals[0,2] = 2020.8334
$rowvals[0,3] = 4471.125
$rowvals[0,4] = 6062.5002
$rowvals[0,5] = -2021.125
$rowvals[0,6] = -10962.5002
$ws.Range("H122:N122").Value = $rowvals

$ws = $wb.Worksheets.Item("BSM")
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 993.36
$rowvals[0,1] = 920.7619
$rowvals[0,2] = 1374.5
$rowvals[0,3] = 920.7619
$rowvals[0,4] = 1374.5
$rowvals[0,5] = -806.7619
$rowvals[0,6] = -1602.5
$ws.Range("H3:N3").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 20471
$rowvals[0,1] = 20471
$rowvals[0,2] = 0
$rowvals[0,3] = 20471
$rowvals[0,4] = 0
$rowvals[0,5] = -20179
$rowvals[0,6] = $null
$ws.Range("H26:N26").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 590.2143
$rowvals[0,1] = 538.4286
$rowvals[0,2] = 642
$rowvals[0,3] = 538.4286
$rowvals[0,4] = 642
$rowvals[0,5] = -313.4286
$rowvals[0,6] = -1092
$ws.Range("H64:N64").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 590.2143
$rowvals[0,1] = 538.4286
$rowvals[0,2] = 642
$rowvals[0,3] = 538.4286
$rowvals[0,4] = 642
$rowvals[0,5] = 241.5714
$rowvals[0,6] = -2202
$ws.Range("H67:N67").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 11598.75
$rowvals[0,1] = 6131.6665
$rowvals[0,2] = 28000
$rowvals[0,3] = 6131.6665
$rowvals[0,4] = 28000
$rowvals[0,5] = -3385.6665
$rowvals[0,6] = -33492
$ws.Range("H96:N96").Value = $rowvals

$ws = $wb.Worksheets.Item("CRP")
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4573.909
$rowvals[0,1] = 0
$rowvals[0,2] = 4573.909
$rowvals[0,3] = 0
$rowvals[0,4] = 4573.909
$rowvals[0,5] = $null
$rowvals[0,6] = -5163.909
$ws.Range("H31:N31").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4573.909
$rowvals[0,1] = 0
$rowvals[0,2] = 4573.909
$rowvals[0,3] = 0
$rowvals[0,4] = 4573.909
$rowvals[0,5] = $null
$rowvals[0,6] = -4977.909
$ws.Range("H34:N34").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 2220073.5
$rowvals[0,1] = 3954461.5
$rowvals[0,2] = 3911.111
$rowvals[0,3] = 3954461.5
$rowvals[0,4] = 3911.111
$rowvals[0,5] = -3954258.5
$rowvals[0,6] = -4317.111
$ws.Range("H58:N58").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 47500
$rowvals[0,1] = 0
$rowvals[0,2] = 47500
$rowvals[0,3] = 0
$rowvals[0,4] = 47500
$rowvals[0,5] = $null
$rowvals[0,6] = -52492
$ws.Range("H92:N92").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 2530.2856
$rowvals[0,1] = 900.3333
$rowvals[0,2] = 3752.75
$rowvals[0,3] = 900.3333
$rowvals[0,4] = 3752.75
$rowvals[0,5] = 846.6667
$rowvals[0,6] = -7246.75
$ws.Range("H105:N105").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 48100
$rowvals[0,1] = 0
$rowvals[0,2] = 48100
$rowvals[0,3] = 0
$rowvals[0,4] = 48100
$rowvals[0,5] = $null
$rowvals[0,6] = -50624
$ws.Range("H106:N106").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 6926.963
$rowvals[0,1] = 2943.158
$rowvals[0,2] = 16388.5
$rowvals[0,3] = 8829.474
$rowvals[0,4] = 49165.5
$rowvals[0,5] = -6379.474
$rowvals[0,6] = -54065.5
$ws.Range("H122:N122").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 2220073.5
$rowvals[0,1] = 3954461.5
$rowvals[0,2] = 3911.111
$rowvals[0,3] = 11863384.5
$rowvals[0,4] = 11733.333
$rowvals[0,5] = -11860834.5
$rowvals[0,6] = -16833.333
$ws.Range("H136:N136").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 40261.555
$rowvals[0,1] = 0
$rowvals[0,2] = 40261.555
$rowvals[0,3] = 0
$rowvals[0,4] = 40261.555
$rowvals[0,5] = $null
$rowvals[0,6] = -50621.555
$ws.Range("H141:N141").Value = $rowvals

$ws = $wb.Worksheets.Item("CUL")
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4007.8386
$rowvals[0,1] = 460.36365
$rowvals[0,2] = 12679.444
$rowvals[0,3] = 1381.09095
$rowvals[0,4] = 38038.33199999999
$rowvals[0,5] = -1269.09095
$rowvals[0,6] = -38262.33199999999
$ws.Range("H5:N5").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 1911.091
$rowvals[0,1] = 1229.1111
$rowvals[0,2] = 4980
$rowvals[0,3] = 3687.3333
$rowvals[0,4] = 14940
$rowvals[0,5] = -3372.3333
$rowvals[0,6] = -15570
$ws.Range("H70:N70").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 1911.091
$rowvals[0,1] = 1229.1111
$rowvals[0,2] = 4980
$rowvals[0,3] = 3687.3333
$rowvals[0,4] = 14940
$rowvals[0,5] = -2595.3333
$rowvals[0,6] = -17124
$ws.Range("H73:N73").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 6800
$rowvals[0,1] = 0
$rowvals[0,2] = 6800
$rowvals[0,3] = 0
$rowvals[0,4] = 20400
$rowvals[0,5] = $null
$rowvals[0,6] = -25642
$ws.Range("H105:N105").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4007.8386
$rowvals[0,1] = 460.36365
$rowvals[0,2] = 12679.444
$rowvals[0,3] = 4143.27285
$rowvals[0,4] = 114114.996
$rowvals[0,5] = -1608.27285
$rowvals[0,6] = -119184.996
$ws.Range("H135:N135").Value = $rowvals

$ws = $wb.Worksheets.Item("GSM")
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 5303.304
$rowvals[0,1] = 4458.25
$rowvals[0,2] = 6225.1816
$rowvals[0,3] = 4458.25
$rowvals[0,4] = 6225.1816
$rowvals[0,5] = -4188.25
$rowvals[0,6] = -6765.1816
$ws.Range("H70:N70").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 5303.304
$rowvals[0,1] = 4458.25
$rowvals[0,2] = 6225.1816
$rowvals[0,3] = 4458.25
$rowvals[0,4] = 6225.1816
$rowvals[0,5] = -3522.25
$rowvals[0,6] = -8097.1816
$ws.Range("H73:N73").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4010.4546
$rowvals[0,1] = 3162.25
$rowvals[0,2] = 5315.385
$rowvals[0,3] = 3162.25
$rowvals[0,4] = 5315.385
$rowvals[0,5] = -1540.25
$rowvals[0,6] = -8559.385
$ws.Range("H102:N102").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4808.1177
$rowvals[0,1] = 4859.6
$rowvals[0,2] = 4422
$rowvals[0,3] = 14578.8
$rowvals[0,4] = 13266
$rowvals[0,5] = -12128.8
$rowvals[0,6] = -18166
$ws.Range("H122:N122").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 2335.4783
$rowvals[0,1] = 1806.0667
$rowvals[0,2] = 3328.125
$rowvals[0,3] = 5418.2001
$rowvals[0,4] = 9984.375
$rowvals[0,5] = -2888.2001
$rowvals[0,6] = -15044.375
$ws.Range("H132:N132").Value = $rowvals

$ws = $wb.Worksheets.Item("LTW")
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4912.9375
$rowvals[0,1] = 3588.5
$rowvals[0,2] = 6237.375
$rowvals[0,3] = 3588.5
$rowvals[0,4] = 6237.375
$rowvals[0,5] = -3476.5
$rowvals[0,6] = -6461.375
$ws.Range("H7:N7").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 11234.482
$rowvals[0,1] = 13774.474
$rowvals[0,2] = 6408.5
$rowvals[0,3] = 13774.474
$rowvals[0,4] = 6408.5
$rowvals[0,5] = -13572.474
$rowvals[0,6] = -6812.5
$ws.Range("H61:N61").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 11234.482
$rowvals[0,1] = 13774.474
$rowvals[0,2] = 6408.5
$rowvals[0,3] = 13774.474
$rowvals[0,4] = 6408.5
$rowvals[0,5] = -11604.474
$rowvals[0,6] = -10748.5
$ws.Range("H113:N113").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 5955.9546
$rowvals[0,1] = 5664.343
$rowvals[0,2] = 7090
$rowvals[0,3] = 16993.029
$rowvals[0,4] = 21270
$rowvals[0,5] = -14543.029
$rowvals[0,6] = -26170
$ws.Range("H122:N122").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 4912.9375
$rowvals[0,1] = 3588.5
$rowvals[0,2] = 6237.375
$rowvals[0,3] = 10765.5
$rowvals[0,4] = 18712.125
$rowvals[0,5] = -8295.5
$rowvals[0,6] = -23652.125
$ws.Range("H126:N126").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 5580.622
$rowvals[0,1] = 3945.6667
$rowvals[0,2] = 8033.0557
$rowvals[0,3] = 11837.0001
$rowvals[0,4] = 24099.1671
$rowvals[0,5] = -9287.000100000001
$rowvals[0,6] = -29199.1671
$ws.Range("H136:N136").Value = $rowvals

$ws = $wb.Worksheets.Item("WVR")
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 3341.75
$rowvals[0,1] = 1886.6875
$rowvals[0,2] = 6251.875
$rowvals[0,3] = 5660.0625
$rowvals[0,4] = 18755.625
$rowvals[0,5] = -3210.0625
$rowvals[0,6] = -23655.625
$ws.Range("H122:N122").Value = $rowvals
$rowvals = New-Object 'object[,]' 1,7
$rowvals[0,0] = 2192.8064
$rowvals[0,1] = 1516.7556
$rowvals[0,2] = 3982.353
$rowvals[0,3] = 3945.6667
$rowvals[0,4] = 11947.059
$rowvals[0,5] = -2000.266799999999
$rowvals[0,6] = -17047.059
$ws.Range("H136:N136").Value = $rowvals
